# Performance_90_PP.xlsx — "Weighting & Scaling update & heatmap"
#
# 1. Update the Scaling sheet's Min/Max values and add two new columns
#    (Optimal / Threshold).
# 2. Make the performance_params_90 sheet the active/selected tab again
#    (it had been on Scaling) and restore its last-used selection.
# 3. Resize column A of the Scaling sheet.

$wb  = $excel.ActiveWorkbook
$wsPerf    = $wb.Worksheets.Item("performance_params_90")
$wsScaling = $wb.Worksheets.Item("Scaling")

# --- Scaling sheet: Min / Max value updates -------------------------------
$wsScaling.Range("B2").Value = 410
$wsScaling.Range("C2").Value = 660

$wsScaling.Range("B3").Value = 600
$wsScaling.Range("C3").Value = 659

$wsScaling.Range("B4").Value = 13
$wsScaling.Range("C4").Value = 18.2

# --- Scaling sheet: new Optimal / Threshold columns ------------------------
$wsScaling.Range("E1").Value = "Optimal"
$wsScaling.Range("F1").Value = "Threshold"

# Match the look of the existing header cells (D1) for the two new ones.
$wsScaling.Range("D1").Copy()
$wsScaling.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Scaling sheet: column A width -----------------------------------------
$wsScaling.Columns.Item(1).ColumnWidth = 10

# --- Sheet selections / active tab -----------------------------------------
$wsScaling.Range("C5").Select()
$wsPerf.Range("B16").Select()
$wsPerf.Activate()
